# Rewrites the title/author block so it matches the Pandoc "title block"
# run-per-word layout:
#   Paragraph 1 ("Thanksgiving Dinner and Other Things", Heading1) becomes
#     a Title-styled paragraph with the text split into one run per word
#     (and a separate run for each inter-word space).
#   Paragraph 2 ("By Dorothy Day", bold) loses the "By " prefix, the bold
#     formatting, and becomes an Authors-styled paragraph, again split
#     into one run per word/space.

$d = $word.ActiveDocument

# --- Paragraph 1: title -----------------------------------------------
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range

$titleXml = '<?xml version="1.0"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:r><w:t xml:space="preserve">Thanksgiving</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
              '<w:r><w:t xml:space="preserve">Dinner</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
              '<w:r><w:t xml:space="preserve">and</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
              '<w:r><w:t xml:space="preserve">Other</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
              '<w:r><w:t xml:space="preserve">Things</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$r1.InsertXML($titleXml)

$p1 = $d.Paragraphs(1)
$p1.Style = "Title"

# --- Paragraph 2: author -------------------------------------------------
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range

$authorXml = '<?xml version="1.0"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
              '<w:r><w:t xml:space="preserve">Day</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$r2.InsertXML($authorXml)

$p2 = $d.Paragraphs(2)
$p2.Style = "Authors"
